# Applies the "BB Alert Result Statistics 7/26" update:
#  - Updates count (B) and mean (D) values for several categories
#  - Updates max (C9) for Roulette_Bullish
#  - Changes the sheet zoom from 210% to 240% and drops the stale selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Rapidfire_Bullish
$ws.Range("B5").Value = 110
$ws.Range("D5").Value = 53.668545454545459

# Row 6 - Repeater_Bearish
$ws.Range("B6").Value = 114
$ws.Range("D6").Value = 60.574649122807052

# Row 7 - Repeater_Bullish
$ws.Range("B7").Value = 396
$ws.Range("D7").Value = 41.434065656565693

# Row 8 - Roulette_Bearish
$ws.Range("B8").Value = 308
$ws.Range("D8").Value = 46.914610389610367

# Row 9 - Roulette_Bullish
$ws.Range("B9").Value = 234
$ws.Range("C9").Value = 416.85
$ws.Range("D9").Value = 41.268547008546989

# Row 10 - Steady_Bullish
$ws.Range("B10").Value = 132
$ws.Range("D10").Value = 38.441666666666663

# Row 11 - Swift_Bullish
$ws.Range("B11").Value = 112
$ws.Range("D11").Value = 39.052321428571418

# Update the view: zoom 210% -> 240%, and clear the old selection (G17) by
# reselecting the top-left cell so the saved sheetView has no stale <selection>.
$ws.Application.ActiveWindow.Zoom = 240
$ws.Range("A1").Select()
